$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAY_2024")

# --- Update the attendance numbers in column E (and related totals) ---
$ws.Range("E10").Value = 25

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 6

$ws.Range("E14").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("E19").Value = 3
$ws.Range("E20").Value = 3
$ws.Range("E21").Value = 3
$ws.Range("E22").Value = 3
$ws.Range("E23").Value = 0

# The %age column divides by 6 (total classes) instead of 3
$ws.Range("G14").Formula = "=(F14/6)*100"
$ws.Range("G15:G23").Formula = "=(F15/6)*100"

# --- Update the view state (selection / scroll position) ---
$ws.Activate()
$ws.Range("E25:G25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 11
